$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to hold a tiny 2-column placeholder table (rows 1-5, with
# several blank rows). The new content is a 7-column x 4-row energy-mix
# table ("Year" plus one column per energy source), so drop the now unused
# trailing row first.
$ws.Rows(5).Delete()

# Helper: write a value into a cell while forcing it to be stored as text
# (shared string), even for the numeric-looking values. This mirrors the
# PDF-extraction source data, where every cell - including figures such as
# "4,251.0" or "2.7" - comes through as plain text rather than a number.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Values are (re)written column by column, top to bottom, matching the
# order the extraction/export step produced the sheet in.
Set-TextValue "A1" "Year"
Set-TextValue "A2" "2019"
Set-TextValue "A3" "2018"
Set-TextValue "A4" "2017"

Set-TextValue "B1" "Wind power"
Set-TextValue "B2" "4,251.0"
Set-TextValue "B3" "3,463.2"
Set-TextValue "B4" "3,683.3"

Set-TextValue "C1" "Biogas"
Set-TextValue "C2" "101.2"
Set-TextValue "C3" "103.2"
Set-TextValue "C4" "82.6"

Set-TextValue "D1" "Biomass"
Set-TextValue "D2" "59.7"
Set-TextValue "D3" "28.8"
Set-TextValue "D4" "12.9"

Set-TextValue "E1" "Photovoltaics"
Set-TextValue "E2" "969.4"
Set-TextValue "E3" "895.2"
Set-TextValue "E4" "1,004.1"

Set-TextValue "F1" "Hydropower"
Set-TextValue "F2" "2.7"
Set-TextValue "F3" "8.3"
Set-TextValue "F4" "8.8"

Set-TextValue "G1" "Total"
Set-TextValue "G2" "5,384.0"
Set-TextValue "G3" "4,498.7"
Set-TextValue "G4" "4,791.7"

# Column A keeps the workbook's original bordered / bold / centered look
# (style index 1 in the source file). The Text-format round trip used by
# Set-TextValue above resets each touched cell back to the Normal style,
# so reapply that formatting now across the whole column in one go.
$colA = $ws.Range("A1:A4")
$colA.Borders.LineStyle = 1
$colA.Font.Bold = $true
$colA.HorizontalAlignment = -4108
$colA.VerticalAlignment = -4160
